$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start the repairs table fresh, then rebuild it with the updated data
$ws.Range("A1:B12").Clear()

# Service ID / Service Type data rows
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "REPAIRS"

$ws.Range("A3").Value = 58
$ws.Range("B3").Value = "HOOD REPAIR"

$ws.Range("A4").Value = 63
$ws.Range("B4").Value = "REPAIR BED FLOOR"

$ws.Range("A5").Value = 74
$ws.Range("B5").Value = "REPAIR TRAILER"

$ws.Range("A6").Value = 98
$ws.Range("B6").Value = "REPAIR DUMP BOX"

# Header row, written last, bolded
$ws.Range("A1").Value = "Service ID"
$ws.Range("B1").Value = "Service Type"

# Restore the thin box border on every cell (Clear wiped formatting)
$ws.Range("A1:B6").Borders.LineStyle = 1
$ws.Range("A1:B1").Font.Bold = $true

# Widen column B slightly to fit the new, longer labels
$ws.Columns.Item(2).ColumnWidth = 18.17

# Move the active selection off the data range (matches the saved UI state)
$ws.Range("G15").Select()
